# Slide 2: subtitle placeholder shape ("副標題 2", shape #2) currently reads
# "While Loop" as a bulleted line under the slide title. Update it to
# "Do While Loop" to match the corrected lesson title.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Clear the existing run/paragraph-end formatting first, then retype the
# text, so the single run ends up holding the corrected text without a
# leftover duplicate end-of-paragraph run.
$tr.Delete()
$tr.Text = "Do While Loop"
